$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.24999533333333
$ws.Range("H2").Value = 99.74998599999999
$ws.Range("I2").Value = 0.5673360890306117
$ws.Range("J2").Value = 0.5673360890306117
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.379101
$ws.Range("N2").Value = 25.137303
$ws.Range("O2").Value = 0.2232365200207407
$ws.Range("P2").Value = 0.2232365200207407
$ws.Range("Q2").Value = 278.6050691475286
$ws.Range("R2").Value = 2507.445622327758
$ws.Range("S2").Value = 0.1266501341973709
$ws.Range("T2").Value = 0.1266501341973709

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.24999533333333
$ws.Range("H3").Value = 99.74998599999999
$ws.Range("I3").Value = 0.5673360890306117
$ws.Range("J3").Value = 0.5673360890306117
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.350178333333334
$ws.Range("N3").Value = 28.050535
$ws.Range("O3").Value = 0.2491080215773342
$ws.Range("P3").Value = 0.2491080215773342
$ws.Range("Q3").Value = 310.8933859491678
$ws.Range("R3").Value = 2798.04047354251
$ws.Range("S3").Value = 0.141327970707838
$ws.Range("T3").Value = 0.141327970707838

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 33.24999533333333
$ws.Range("H4").Value = 99.74998599999999
$ws.Range("I4").Value = 0.5673360890306117
$ws.Range("J4").Value = 0.5673360890306117
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.80535433333333
$ws.Range("N4").Value = 59.416063
$ws.Range("O4").Value = 0.5276554584019252
$ws.Range("P4").Value = 0.5276554584019252
$ws.Range("Q4").Value = 658.5279391583464
$ws.Range("R4").Value = 5926.751452425117
$ws.Range("S4").Value = 0.2993579841254029
$ws.Range("T4").Value = 0.2993579841254029

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.30243966666667
$ws.Range("H5").Value = 69.907319
$ws.Range("I5").Value = 0.3976035140102714
$ws.Range("J5").Value = 0.3976035140102714
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.379101
$ws.Range("N5").Value = 25.137303
$ws.Range("O5").Value = 0.2232365200207407
$ws.Range("P5").Value = 0.2232365200207407
$ws.Range("Q5").Value = 195.2534955134064
$ws.Range("R5").Value = 1757.281459620657
$ws.Range("S5").Value = 0.0887596248156708
$ws.Range("T5").Value = 0.08875962481567079

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.30243966666667
$ws.Range("H6").Value = 69.907319
$ws.Range("I6").Value = 0.3976035140102714
$ws.Range("J6").Value = 0.3976035140102714
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.350178333333334
$ws.Range("N6").Value = 28.050535
$ws.Range("O6").Value = 0.2491080215773342
$ws.Range("P6").Value = 0.2491080215773342
$ws.Range("Q6").Value = 217.8819664850739
$ws.Range("R6").Value = 1960.937698365665
$ws.Range("S6").Value = 0.0990462247472946
$ws.Range("T6").Value = 0.09904622474729458

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.30243966666667
$ws.Range("H7").Value = 69.907319
$ws.Range("I7").Value = 0.3976035140102714
$ws.Range("J7").Value = 0.3976035140102714
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.80535433333333
$ws.Range("N7").Value = 59.416063
$ws.Range("O7").Value = 0.5276554584019252
$ws.Range("P7").Value = 0.5276554584019252
$ws.Range("Q7").Value = 461.5130744294553
$ws.Range("R7").Value = 4153.617669865097
$ws.Range("S7").Value = 0.209797664447306
$ws.Range("T7").Value = 0.209797664447306

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.054792666666667
$ws.Range("H8").Value = 6.164378
$ws.Range("I8").Value = 0.03506039695911681
$ws.Range("J8").Value = 0.03506039695911681
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.379101
$ws.Range("N8").Value = 25.137303
$ws.Range("O8").Value = 0.2232365200207407
$ws.Range("P8").Value = 0.2232365200207407
$ws.Range("Q8").Value = 17.21731528805934
$ws.Range("R8").Value = 154.955837592534
$ws.Range("S8").Value = 0.007826761007698996
$ws.Range("T8").Value = 0.007826761007698994

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.054792666666667
$ws.Range("H9").Value = 6.164378
$ws.Range("I9").Value = 0.03506039695911681
$ws.Range("J9").Value = 0.03506039695911681
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.350178333333334
$ws.Range("N9").Value = 28.050535
$ws.Range("O9").Value = 0.2491080215773342
$ws.Range("P9").Value = 0.2491080215773342
$ws.Range("Q9").Value = 19.21267787135889
$ws.Range("R9").Value = 172.91410084223
$ws.Range("S9").Value = 0.008733826122201574
$ws.Range("T9").Value = 0.00873382612220157

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.054792666666667
$ws.Range("H10").Value = 6.164378
$ws.Range("I10").Value = 0.03506039695911681
$ws.Range("J10").Value = 0.03506039695911681
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.80535433333333
$ws.Range("N10").Value = 59.416063
$ws.Range("O10").Value = 0.5276554584019252
$ws.Range("P10").Value = 0.5276554584019252
$ws.Range("Q10").Value = 40.69589684486822
$ws.Range("R10").Value = 366.263071603814
$ws.Range("S10").Value = 0.01849980982921625
$ws.Range("T10").Value = 0.01849980982921624

